# Adapt column header formatting to respective input file names:
#   *_old -> *_FV2310
#   *_new -> *_FV2404
# Then turn the data range into an Excel Table ("Table1") and freeze the
# header row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldHeaders = @(
    "Segmentname_old",
    "Segmentgruppe_old",
    "Segment_old",
    "Datenelement_old",
    "Segment ID_old",
    "Code_old",
    "Qualifier_old",
    "Beschreibung_old",
    "Bedingungsausdruck_old",
    "Bedingung_old"
)

$newHeaders = @(
    "Segmentname_new",
    "Segmentgruppe_new",
    "Segment_new",
    "Datenelement_new",
    "Segment ID_new",
    "Code_new",
    "Qualifier_new",
    "Beschreibung_new",
    "Bedingungsausdruck_new",
    "Bedingung_new"
)

# Columns A..J (1..10) carry the "_old" suffix -> rename to "_FV2310"
for ($i = 0; $i -lt $oldHeaders.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 1)
    $suffixFree = $oldHeaders[$i].Substring(0, $oldHeaders[$i].Length - 4)
    $cell.Value = "$suffixFree" + "_FV2310"
}

# Column K (11) is "diff" and stays untouched.

# Columns L..U (12..21) carry the "_new" suffix -> rename to "_FV2404"
for ($i = 0; $i -lt $newHeaders.Count; $i++) {
    $cell = $ws.Cells.Item(1, $i + 12)
    $suffixFree = $newHeaders[$i].Substring(0, $newHeaders[$i].Length - 4)
    $cell.Value = "$suffixFree" + "_FV2404"
}

# Convert the used data range into an Excel Table so the header row also
# gets the filter dropdowns, and freeze the header row.
$dataRange = $ws.Range("A1:U55")
$table = $ws.ListObjects.Add(1, $dataRange, 0, 1)
$table.Name = "Table1"

$ws.Activate()
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
